$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("planificación")
$ws.Range("O6").Value = -1
